$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C ("Förändrad") for rows 2-9 is being updated from serial date
# 45233 (2023-11-03) to 45243 (2023-11-13).
$ws.Range("C2:C9").Value = 45243
